$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 17.08608867836142)
    3 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    4 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
    5 = @(0.6545652718822623, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.716211508195562)
    6 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
